# Apply the "Handles float input without breaking stuff" edit:
#  - Recomputed marksheet summary numbers (rows 10-12)
#  - Collapsed the 3 repeated "Student Ans / Correct Ans" attempt blocks
#    (columns A/B, D/E, G/H) down to a single populated block, filling in
#    each student's chosen answer (colored by correct/incorrect/not-attempted)
#    and clearing out the now-unused G/H block entirely, and most of the D/E
#    block (kept only for rows 16-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "63/112"

# ---- Remove the third "Student Ans / Correct Ans" block (columns G/H) ----
$ws.Range("G15:H40").Clear()

# ---- Remove the second block (columns D/E) except for rows 16-18 ----
$ws.Range("D19:E40").Clear()

# ---- Fill in student answers for the first block (column A) ----
# Correct answers (style/value copied from the graded "correctStyle")
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A16").Value = "Option A"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"

$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"

$ws.Range("A20").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A23").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"

$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"

$ws.Range("A27").Style = "incorrectStyle"
$ws.Range("A27").Value = "Option D"

$ws.Range("A28").Style = "correctStyle"
$ws.Range("A28").Value = "Option D"

$ws.Range("A29").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"

$ws.Range("A30").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"

$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"

$ws.Range("A34").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"

$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").Value = "Option A"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"

# Rows 17, 22, 24, 26, 31, 32, 35, 37, 38, 40 remain "not attempted"
# (already blank with normalStyle), so nothing to do for those.

# ---- Fill in the kept D column answers for rows 16-18 ----
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"

# D17 remains "not attempted" (blank, normalStyle) - nothing to do.
